$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 92, shifting existing rows 92-95 down to 93-96.
$ws.Rows.Item(92).Insert()

# Copy formatting (number format / style) from the row that is now 93 (previously row 92)
# into the newly inserted row 92, so the date cell keeps its date style.
# Restrict the copy to the used columns (A:R) so we don't materialize formatting
# across the entire 16384-column row.
$ws.Range("A93:R93").Copy()
$ws.Range("A92:R92").PasteSpecial(-4122) # xlPasteFormats

# Fill in the values for the new row 92 (matches the existing pattern for this data set).
$ws.Cells.Item(92, 1).Value = 8
$ws.Cells.Item(92, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 44516
$ws.Cells.Item(92, 5).Value = 4
$ws.Cells.Item(92, 6).Value = 100112044
$ws.Cells.Item(92, 7).Value = "Perejil"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 3100
$ws.Cells.Item(92, 11).Value = 1300
$ws.Cells.Item(92, 12).Value = 1500
$ws.Cells.Item(92, 13).Value = 1400
$ws.Cells.Item(92, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(92, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(92, 16).Value = 933
$ws.Cells.Item(92, 17).Value = 1.5
$ws.Cells.Item(92, 18).Value = "Hortaliza"
